$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 912.82355
$ws.Range("I28").Value = 868.2
$ws.Range("J28").Value = 976.5714
$ws.Range("K28").Value = 868.2
$ws.Range("L28").Value = 976.5714
$ws.Range("M28").Value = -383.2
$ws.Range("N28").Value = -1946.5714

$ws.Range("H70").Value = 829
$ws.Range("J70").Value = 1275.7778
$ws.Range("L70").Value = 3827.3334
$ws.Range("N70").Value = -4367.3334

$ws.Range("H73").Value = 829
$ws.Range("J73").Value = 1275.7778
$ws.Range("L73").Value = 3827.3334
$ws.Range("N73").Value = -5699.3334

$ws.Range("H86").Value = 5039.3125
$ws.Range("I86").Value = 3916.7778
$ws.Range("J86").Value = 6482.5713
$ws.Range("K86").Value = 3916.7778
$ws.Range("L86").Value = 6482.5713
$ws.Range("M86").Value = -2793.7778
$ws.Range("N86").Value = -8728.5713

$ws.Range("H89").Value = 5039.3125
$ws.Range("I89").Value = 3916.7778
$ws.Range("J89").Value = 6482.5713
$ws.Range("K89").Value = 19583.889
$ws.Range("L89").Value = 32412.8565
$ws.Range("M89").Value = -13967.889
$ws.Range("N89").Value = -43644.85649999999

$ws.Range("H95").Value = 32162
$ws.Range("J95").Value = 32162
$ws.Range("L95").Value = 32162
$ws.Range("N95").Value = -37654

$ws.Range("N96").ClearContents()
$ws.Range("H96").Value = 689
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0

$ws.Range("H97").Value = 15000
$ws.Range("J97").Value = 15000
$ws.Range("L97").Value = 45000
$ws.Range("N97").Value = -45992

$ws.Range("H141").Value = 4769.7744
$ws.Range("I141").Value = 3949
$ws.Range("K141").Value = 11847
$ws.Range("M141").Value = -6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N27").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0

$ws.Range("H45").Value = 3318.4443
$ws.Range("I45").Value = 3111.0833
$ws.Range("J45").Value = 3733.1667
$ws.Range("K45").Value = 3111.0833
$ws.Range("L45").Value = 3733.1667
$ws.Range("M45").Value = -2734.0833
$ws.Range("N45").Value = -4487.1667

$ws.Range("H61").Value = 3643.724
$ws.Range("I61").Value = 3907.04
$ws.Range("K61").Value = 3907.04
$ws.Range("M61").Value = -3695.04

$ws.Range("H74").Value = 1733.8
$ws.Range("I74").Value = 1861.0588
$ws.Range("K74").Value = 1861.0588
$ws.Range("M74").Value = -987.0588

$ws.Range("H76").Value = 84833.336
$ws.Range("J76").Value = 84833.336
$ws.Range("L76").Value = 84833.336
$ws.Range("N76").Value = -85509.336

$ws.Range("H77").Value = 1733.8
$ws.Range("I77").Value = 1861.0588
$ws.Range("K77").Value = 9305.294
$ws.Range("M77").Value = -4937.294

$ws.Range("H79").Value = 84833.336
$ws.Range("J79").Value = 84833.336
$ws.Range("L79").Value = 84833.336
$ws.Range("N79").Value = -87173.336

$ws.Range("H132").Value = 3250.4194
$ws.Range("I132").Value = 3299.077
$ws.Range("J132").Value = 2997.4
$ws.Range("K132").Value = 9897.231
$ws.Range("L132").Value = 8992.200000000001
$ws.Range("M132").Value = -7367.231
$ws.Range("N132").Value = -14052.2

$ws.Range("H136").Value = 3643.724
$ws.Range("I136").Value = 3907.04
$ws.Range("K136").Value = 11721.12
$ws.Range("M136").Value = -9171.119999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 27779908
$ws.Range("I86").Value = 45455764
$ws.Range("J86").Value = 3561.4285
$ws.Range("K86").Value = 45455764
$ws.Range("L86").Value = 3561.4285
$ws.Range("M86").Value = -45454641
$ws.Range("N86").Value = -5807.4285

$ws.Range("H89").Value = 27779908
$ws.Range("I89").Value = 45455764
$ws.Range("J89").Value = 3561.4285
$ws.Range("K89").Value = 227278820
$ws.Range("L89").Value = 17807.1425
$ws.Range("M89").Value = -227273204
$ws.Range("N89").Value = -29039.1425

$ws.Range("H96").Value = 22665.334
$ws.Range("J96").Value = 24855.572
$ws.Range("L96").Value = 24855.572
$ws.Range("N96").Value = -30347.572

$ws.Range("H99").Value = 53001.25
$ws.Range("I99").Value = 101127.5
$ws.Range("J99").Value = 4875
$ws.Range("K99").Value = 101127.5
$ws.Range("L99").Value = 4875
$ws.Range("M99").Value = -99629.5
$ws.Range("N99").Value = -7871

$ws.Range("H105").Value = 2816.3333
$ws.Range("I105").Value = 1254.091
$ws.Range("K105").Value = 1254.091
$ws.Range("M105").Value = 492.9090000000001

$ws.Range("H134").Value = 2431
$ws.Range("I134").Value = 1890.3636
$ws.Range("J134").Value = 2780.8235
$ws.Range("K134").Value = 5671.0908
$ws.Range("L134").Value = 8342.470499999999
$ws.Range("M134").Value = -3136.0908
$ws.Range("N134").Value = -13412.4705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 612.7143
$ws.Range("J22").Value = 666.3333
$ws.Range("L22").Value = 666.3333
$ws.Range("N22").Value = -1366.3333

$ws.Range("H31").Value = 1618.6279
$ws.Range("I31").Value = 1464.675
$ws.Range("K31").Value = 1464.675
$ws.Range("M31").Value = -1169.675

$ws.Range("H34").Value = 1618.6279
$ws.Range("I34").Value = 1464.675
$ws.Range("K34").Value = 1464.675
$ws.Range("M34").Value = -1262.675

$ws.Range("H41").Value = 28061.5
$ws.Range("I41").Value = 3750
$ws.Range("J41").Value = 44269.168
$ws.Range("K41").Value = 3750
$ws.Range("L41").Value = 44269.168
$ws.Range("M41").Value = -3322
$ws.Range("N41").Value = -45125.168

$ws.Range("H58").Value = 1778.619
$ws.Range("I58").Value = 1488.2941
$ws.Range("K58").Value = 1488.2941
$ws.Range("M58").Value = -1285.2941

$ws.Range("H105").Value = 983.3
$ws.Range("I105").Value = 971.6
$ws.Range("K105").Value = 971.6
$ws.Range("M105").Value = 775.4

$ws.Range("H132").Value = 2111.6072
$ws.Range("I132").Value = 1558.2609
$ws.Range("K132").Value = 4674.7827
$ws.Range("M132").Value = -2144.7827

$ws.Range("H134").Value = 1285.5714
$ws.Range("I134").Value = 1285.5714
$ws.Range("K134").Value = 3856.7142
$ws.Range("M134").Value = -1321.7142

$ws.Range("H136").Value = 1778.619
$ws.Range("I136").Value = 1488.2941
$ws.Range("K136").Value = 4464.8823
$ws.Range("M136").Value = -1914.8823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6846.857
$ws.Range("I80").Value = 4497.143
$ws.Range("J80").Value = 9196.571
$ws.Range("K80").Value = 4497.143
$ws.Range("L80").Value = 9196.571
$ws.Range("M80").Value = -3499.143
$ws.Range("N80").Value = -11192.571

$ws.Range("H83").Value = 6846.857
$ws.Range("I83").Value = 4497.143
$ws.Range("J83").Value = 9196.571
$ws.Range("K83").Value = 22485.715
$ws.Range("L83").Value = 45982.855
$ws.Range("M83").Value = -17493.715
$ws.Range("N83").Value = -55966.855

$ws.Range("H99").Value = 21976.75
$ws.Range("I99").Value = 8263.799999999999
$ws.Range("K99").Value = 8263.799999999999
$ws.Range("M99").Value = -6017.799999999999

$ws.Range("H107").Value = 434.8421
$ws.Range("I107").Value = 365.92307
$ws.Range("J107").Value = 584.1667
$ws.Range("K107").Value = 365.92307
$ws.Range("L107").Value = 584.1667
$ws.Range("M107").Value = 1554.07693
$ws.Range("N107").Value = -4424.1667

$ws.Range("H132").Value = 7450.59
$ws.Range("I132").Value = 6587.8
$ws.Range("K132").Value = 19763.4
$ws.Range("M132").Value = -17233.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 25000322
$ws.Range("I16").Value = 27778112
$ws.Range("J16").Value = 198
$ws.Range("K16").Value = 27778112
$ws.Range("L16").Value = 198
$ws.Range("M16").Value = -27777942
$ws.Range("N16").Value = -538

$ws.Range("H40").Value = 1749.5
$ws.Range("I40").Value = 1749.5
$ws.Range("K40").Value = 1749.5
$ws.Range("M40").Value = -1613.5

$ws.Range("H55").Value = 537.087
$ws.Range("I55").Value = 538.9
$ws.Range("K55").Value = 538.9
$ws.Range("M55").Value = -365.9

$ws.Range("H100").Value = 6031.2856
$ws.Range("I100").Value = 4761
$ws.Range("J100").Value = 6539.4
$ws.Range("K100").Value = 4761
$ws.Range("L100").Value = 6539.4
$ws.Range("M100").Value = -4220
$ws.Range("N100").Value = -7621.4

$ws.Range("H132").Value = 3122.56
$ws.Range("I132").Value = 2141.4707
$ws.Range("K132").Value = 6424.4121
$ws.Range("M132").Value = -3894.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 30000
$ws.Range("I26").Value = 30000
$ws.Range("K26").Value = 30000
$ws.Range("M26").Value = -29707

$ws.Range("M43").ClearContents()
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20298

$ws.Range("H51").Value = 34284.715
$ws.Range("I51").Value = 20000
$ws.Range("K51").Value = 20000
$ws.Range("M51").Value = -19490

$ws.Range("H126").Value = 4591.3477
$ws.Range("I126").Value = 4873.4287
$ws.Range("K126").Value = 14620.2861
$ws.Range("M126").Value = -12150.2861

$ws.Range("H132").Value = 3044.9092
$ws.Range("I132").Value = 2780.9092
$ws.Range("J132").Value = 4100.909
$ws.Range("K132").Value = 8342.7276
$ws.Range("L132").Value = 12302.727
$ws.Range("M132").Value = -5812.7276
$ws.Range("N132").Value = -17362.727
